# "Generate Report for handoff"
#
# A new handoff was generated for the "131630ab-b982-4819-9e5e-f4b751a73c51"
# file (row 4 of both the zh-cn and de-de report sheets). The "Latest
# Handoff File" (column C) stays the same file name, but the "Latest
# Handoff Datetime" (column D) is refreshed to record the new handoff
# timestamp.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-15 02:57:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-15 02:57:54"
